$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "문과 데이터 사이언스 학과?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ds-for-liberal-arts/#utm_source=rss&utm_medium=rss&utm_campaign=ds-for-liberal-arts"

$ws.Range("D23").Value = "교보문고 컴퓨터공학분야 국내도서 베스트셀러 최근1주일 기준 (2022년 7월 13일)"
$ws.Range("E23").Value = "https://theonly1.tistory.com/3003"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D27").Value = "AWS Inferentia 를 이용한 모델 서빙 비용 최적화: 모델 서버 비용 2배 줄이기 1탄"
$ws.Range("E27").Value = "https://blog.pingpong.us/aws-inferentia/"

$ws.Range("D51").Value = "[MySQL] root 계정 비밀번호 변경 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MySQL-root-%EA%B3%84%EC%A0%95-%EB%B9%84%EB%B0%80%EB%B2%88%ED%98%B8-%EB%B3%80%EA%B2%BD-%EB%B0%A9%EB%B2%95"
